$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= source data (was row 9)
$ws.Range("A2").Value = 111883983
$ws.Range("B2").Value = 90332
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = 'Svavelriska'
$ws.Range("G2").Value = 'Lactarius scrobiculatus'
$ws.Range("H2").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P2").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q2").Value = 575058.3527020445
$ws.Range("R2").Value = 6703446.206921679
$ws.Range("AW2").Value = 'Patric Engfeldt'
$ws.Range("AX2").Value = 'Patric Engfeldt'

# Row 3 <= source data (was row 6)
$ws.Range("A3").Value = 111896640
$ws.Range("B3").Value = 90332
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 4769
$ws.Range("F3").Value = 'Svavelriska'
$ws.Range("G3").Value = 'Lactarius scrobiculatus'
$ws.Range("H3").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P3").Value = 'Kratte masugn, Gstr'
$ws.Range("Q3").Value = 575025.3556637274
$ws.Range("R3").Value = 6703369.042946251
$ws.Range("AW3").Value = 'Philipp Weiss'
$ws.Range("AX3").Value = 'Philipp Weiss'

# Row 4 <= source data (was row 2)
$ws.Range("A4").Value = 111896654
$ws.Range("B4").Value = 89183
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 3215
$ws.Range("F4").Value = 'Rödgul trumpetsvamp'
$ws.Range("G4").Value = 'Craterellus lutescens'
$ws.Range("H4").Value = '(Fr.) Fr.'
$ws.Range("P4").Value = 'Kratte masugn, Gstr'
$ws.Range("Q4").Value = 575072.6962527435
$ws.Range("R4").Value = 6703421.833381963
$ws.Range("AW4").Value = 'Philipp Weiss'
$ws.Range("AX4").Value = 'Philipp Weiss'

# Row 5 <= source data (was row 4)
$ws.Range("A5").Value = 111896642
$ws.Range("B5").Value = 90332
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 4769
$ws.Range("F5").Value = 'Svavelriska'
$ws.Range("G5").Value = 'Lactarius scrobiculatus'
$ws.Range("H5").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P5").Value = 'Kratte masugn, Gstr'
$ws.Range("Q5").Value = 575014.1091647458
$ws.Range("R5").Value = 6703387.066676207
$ws.Range("AW5").Value = 'Philipp Weiss'
$ws.Range("AX5").Value = 'Philipp Weiss'

# Row 6 <= source data (was row 15)
$ws.Range("A6").Value = 111884133
$ws.Range("B6").Value = 88899
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 3286
$ws.Range("F6").Value = 'Flattoppad klubbsvamp'
$ws.Range("G6").Value = 'Clavariadelphus truncatus'
$ws.Range("H6").Value = '(Quél.) Donk'
$ws.Range("P6").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q6").Value = 575059.034285416
$ws.Range("R6").Value = 6703389.477814267
$ws.Range("AW6").Value = 'Patric Engfeldt'
$ws.Range("AX6").Value = 'Patric Engfeldt'

# Row 7 <= source data (was row 21)
$ws.Range("A7").Value = 111896644
$ws.Range("B7").Value = 90332
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 4769
$ws.Range("F7").Value = 'Svavelriska'
$ws.Range("G7").Value = 'Lactarius scrobiculatus'
$ws.Range("H7").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P7").Value = 'Kratte masugn, Gstr'
$ws.Range("Q7").Value = 575036.4083237475
$ws.Range("R7").Value = 6703431.936489306
$ws.Range("AW7").Value = 'Philipp Weiss'
$ws.Range("AX7").Value = 'Philipp Weiss'

# Row 8 <= source data (was row 22)
$ws.Range("A8").Value = 111896639
$ws.Range("B8").Value = 90332
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 4769
$ws.Range("F8").Value = 'Svavelriska'
$ws.Range("G8").Value = 'Lactarius scrobiculatus'
$ws.Range("H8").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P8").Value = 'Kratte masugn, Gstr'
$ws.Range("Q8").Value = 575089.384229039
$ws.Range("R8").Value = 6703379.745088123
$ws.Range("AW8").Value = 'Philipp Weiss'
$ws.Range("AX8").Value = 'Philipp Weiss'

# Row 9 <= source data (was row 17)
$ws.Range("A9").Value = 111896641
$ws.Range("B9").Value = 90332
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 4769
$ws.Range("F9").Value = 'Svavelriska'
$ws.Range("G9").Value = 'Lactarius scrobiculatus'
$ws.Range("H9").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P9").Value = 'Kratte masugn, Gstr'
$ws.Range("Q9").Value = 575021.3626164712
$ws.Range("R9").Value = 6703370.933926445
$ws.Range("AW9").Value = 'Philipp Weiss'
$ws.Range("AX9").Value = 'Philipp Weiss'

# Row 10 <= source data (was row 7)
$ws.Range("A10").Value = 111896637
$ws.Range("B10").Value = 90332
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 4769
$ws.Range("F10").Value = 'Svavelriska'
$ws.Range("G10").Value = 'Lactarius scrobiculatus'
$ws.Range("H10").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P10").Value = 'Kratte masugn, Gstr'
$ws.Range("Q10").Value = 575088.0587098968
$ws.Range("R10").Value = 6703396.00058554
$ws.Range("AW10").Value = 'Philipp Weiss'
$ws.Range("AX10").Value = 'Philipp Weiss'

# Row 11 <= source data (was row 10)
$ws.Range("A11").Value = 111896653
$ws.Range("B11").Value = 89183
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 3215
$ws.Range("F11").Value = 'Rödgul trumpetsvamp'
$ws.Range("G11").Value = 'Craterellus lutescens'
$ws.Range("H11").Value = '(Fr.) Fr.'
$ws.Range("P11").Value = 'Kratte masugn, Gstr'
$ws.Range("Q11").Value = 575075.050630242
$ws.Range("R11").Value = 6703403.625642136
$ws.Range("AW11").Value = 'Philipp Weiss'
$ws.Range("AX11").Value = 'Philipp Weiss'

# Row 12 <= source data (was row 18)
$ws.Range("A12").Value = 111896652
$ws.Range("B12").Value = 89183
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 3215
$ws.Range("F12").Value = 'Rödgul trumpetsvamp'
$ws.Range("G12").Value = 'Craterellus lutescens'
$ws.Range("H12").Value = '(Fr.) Fr.'
$ws.Range("P12").Value = 'Kratte masugn, Gstr'
$ws.Range("Q12").Value = 575066.556649723
$ws.Range("R12").Value = 6703455.751857814
$ws.Range("AW12").Value = 'Philipp Weiss'
$ws.Range("AX12").Value = 'Philipp Weiss'

# Row 13 <= source data (was row 19)
$ws.Range("A13").Value = 111896655
$ws.Range("B13").Value = 89183
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 3215
$ws.Range("F13").Value = 'Rödgul trumpetsvamp'
$ws.Range("G13").Value = 'Craterellus lutescens'
$ws.Range("H13").Value = '(Fr.) Fr.'
$ws.Range("P13").Value = 'Kratte masugn, Gstr'
$ws.Range("Q13").Value = 575104.6742508161
$ws.Range("R13").Value = 6703428.910891063
$ws.Range("AW13").Value = 'Philipp Weiss'
$ws.Range("AX13").Value = 'Philipp Weiss'

# Row 15 <= source data (was row 11)
$ws.Range("A15").Value = 111896633
$ws.Range("B15").Value = 90332
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 4769
$ws.Range("F15").Value = 'Svavelriska'
$ws.Range("G15").Value = 'Lactarius scrobiculatus'
$ws.Range("H15").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P15").Value = 'Kratte masugn, Gstr'
$ws.Range("Q15").Value = 575100.4050603262
$ws.Range("R15").Value = 6703444.118284944
$ws.Range("AW15").Value = 'Philipp Weiss'
$ws.Range("AX15").Value = 'Philipp Weiss'

# Row 17 <= source data (was row 5)
$ws.Range("A17").Value = 111896634
$ws.Range("B17").Value = 90332
$ws.Range("D17").Value = 'LC'
$ws.Range("E17").Value = 4769
$ws.Range("F17").Value = 'Svavelriska'
$ws.Range("G17").Value = 'Lactarius scrobiculatus'
$ws.Range("H17").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P17").Value = 'Kratte masugn, Gstr'
$ws.Range("Q17").Value = 575048.3395925189
$ws.Range("R17").Value = 6703452.413791304
$ws.Range("AW17").Value = 'Philipp Weiss'
$ws.Range("AX17").Value = 'Philipp Weiss'

# Row 18 <= source data (was row 8)
$ws.Range("A18").Value = 111884471
$ws.Range("B18").Value = 88899
$ws.Range("D18").Value = 'NT'
$ws.Range("E18").Value = 3286
$ws.Range("F18").Value = 'Flattoppad klubbsvamp'
$ws.Range("G18").Value = 'Clavariadelphus truncatus'
$ws.Range("H18").Value = '(Quél.) Donk'
$ws.Range("P18").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q18").Value = 575020.8210917887
$ws.Range("R18").Value = 6703397.074168184
$ws.Range("AW18").Value = 'Patric Engfeldt'
$ws.Range("AX18").Value = 'Patric Engfeldt'

# Row 19 <= source data (was row 12)
$ws.Range("A19").Value = 111896635
$ws.Range("B19").Value = 90332
$ws.Range("D19").Value = 'LC'
$ws.Range("E19").Value = 4769
$ws.Range("F19").Value = 'Svavelriska'
$ws.Range("G19").Value = 'Lactarius scrobiculatus'
$ws.Range("H19").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P19").Value = 'Kratte masugn, Gstr'
$ws.Range("Q19").Value = 575037.2974304935
$ws.Range("R19").Value = 6703389.027347369
$ws.Range("AW19").Value = 'Philipp Weiss'
$ws.Range("AX19").Value = 'Philipp Weiss'

# Row 21 <= source data (was row 3)
$ws.Range("A21").Value = 111896636
$ws.Range("B21").Value = 90332
$ws.Range("D21").Value = 'LC'
$ws.Range("E21").Value = 4769
$ws.Range("F21").Value = 'Svavelriska'
$ws.Range("G21").Value = 'Lactarius scrobiculatus'
$ws.Range("H21").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P21").Value = 'Kratte masugn, Gstr'
$ws.Range("Q21").Value = 575108.85141061
$ws.Range("R21").Value = 6703418.142308297
$ws.Range("AW21").Value = 'Philipp Weiss'
$ws.Range("AX21").Value = 'Philipp Weiss'

# Row 22 <= source data (was row 13)
$ws.Range("A22").Value = 111896643
$ws.Range("B22").Value = 90332
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 4769
$ws.Range("F22").Value = 'Svavelriska'
$ws.Range("G22").Value = 'Lactarius scrobiculatus'
$ws.Range("H22").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P22").Value = 'Kratte masugn, Gstr'
$ws.Range("Q22").Value = 575038.7114136803
$ws.Range("R22").Value = 6703416.194821274
$ws.Range("AW22").Value = 'Philipp Weiss'
$ws.Range("AX22").Value = 'Philipp Weiss'

# Clear K column for rows that no longer carry the extra Alder-Stadium marker
$ws.Range("K8").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("K15").ClearContents()

# Rows that gain the (empty) Alder-Stadium marker cell
$ws.Range("K2").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("K18").Value = ""
